$wb = $excel.ActiveWorkbook

# Insert new "Elf" worksheet between "Lifespan" and "Multiplier"
$lifespan = $wb.Worksheets.Item("Lifespan")
$elf = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lifespan)
$elf.Name = "Elf"

# Header row
$elf.Range("A1").Value = "Age"
$elf.Range("B1").Value = "Equivalent"

# Age column: 18, 36, 54 ... 360 (18 * row offset)
for ($i = 0; $i -lt 20; $i++) {
    $r = $i + 2
    $age = ($i + 1) * 18
    $elf.Cells.Item($r, 1).Value = $age
}

# Equivalent column: ROUNDDOWN(Age/4, 0) - B2 standalone, B3:B21 filled as a shared formula
$elf.Range("B2").Formula = "=ROUNDDOWN(A2/4, 0)"
$elf.Range("B3:B21").Formula = "=ROUNDDOWN(A3/4, 0)"

# Match the saved selection/active cell and make Elf the active sheet/tab
$elf.Range("E20").Select() | Out-Null
$elf.Activate() | Out-Null
